$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the added line contingencies (line7, line8);
# this physically pushes the existing extr1..extr8 rows from 8:15 down to
# 10:17, carrying their C/D/E values with them.
$ws.Rows("8:9").Insert()

# Insert() drops the thin border on column A for the freshly inserted
# rows - restore the header/index formatting used throughout column A.
$ws.Range("A8:A9").Borders.LineStyle = 1
$ws.Range("A8:A9").Borders.Weight = 2
$ws.Range("A8:A9").Font.Bold = $true
$ws.Range("A8:A9").HorizontalAlignment = -4108
$ws.Range("A8:A9").VerticalAlignment = -4160

# New row 8: line7
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row 9: line8
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# The index column is a running count (row-2); renumber it for the
# shifted-down extr1..extr8 rows.
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# in_service flips on a few of the shifted-down extr rows
$ws.Range("E10").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("E13").Value = $false
$ws.Range("E16").Value = $false
